# "BOM changes to reflect the now missing 2 buttons"
# Two pushbutton rows (S4 at row 12, S5 at row 13) are no longer part of
# the board, so their BOM rows are removed. Deleting the rows (rather than
# just clearing their contents) shifts the remaining rows (U$3 OLED, U$31
# AAA3 enclosure) up from 14/15 to 12/13, and the sheet's used range shrinks
# from A1:H15 to A1:H13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two BOM rows for the removed buttons (S4, S5) - this shifts
# the rows below them upward.
$ws.Range("A12:A13").EntireRow.Delete()

# Leave the selection on the first of the rows that moved up into view.
$ws.Range("A12").Select()
